# Insert a new data row at row 244 (shifts existing rows 244:302 down to 245:303)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(244).Insert()

# Populate the newly-inserted row 244 with the new record
$ws.Range("A244").Value2 = 4
$ws.Range("B244").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C244").Value2 = "Los Lagos"
$ws.Range("D244").Value2 = 44943
$ws.Range("E244").Value2 = 10
$ws.Range("F244").Value2 = "Fruta"
$ws.Range("G244").Value2 = 100108
$ws.Range("H244").Value2 = "Tropicales y subtropicales"
$ws.Range("I244").Value2 = 100108002
$ws.Range("J244").Value2 = "Mango"
$ws.Range("K244").Value2 = "Sin especificar"
$ws.Range("L244").Value2 = "Primera"
$ws.Range("M244").Value2 = 200
$ws.Range("N244").Value2 = 7500
$ws.Range("O244").Value2 = 8000
$ws.Range("P244").Value2 = 7750
$ws.Range("Q244").Value2 = "$/bandeja 4 kilos"
$ws.Range("R244").Value2 = "Brasil"
$ws.Range("S244").Value2 = 1938
$ws.Range("T244").Value2 = 4
